$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "('Bear', ['Token Creature — Bear', '2/2'])"
$ws.Range("A3").Value = "('Beast', ['Token Creature — Beast', '4/4'])"
$ws.Range("A4").Value = "('Bird', ['Token Creature — Bird', 'Flying', '1/1'])"
$ws.Range("A5").Value = "('Elephant', ['Token Creature — Elephant', '3/3'])"
$ws.Range("A6").Value = "('Goblin Soldier', ['Token Creature — Goblin Soldier', '1/1'])"
$ws.Range("A7").Value = "('Saproling', ['Token Creature — Saproling', '1/1'])"
$ws.Range("A8").Value = "('Spirit', ['Token Creature — Spirit', 'Flying', '1/1'])"
$ws.Range("A9").Value = "('Wasteland', ['Land', '{T}: Add {C}.', '{T}, Sacrifice Wasteland: Destroy target nonbasic land.'])"

$ws.Range("A10:A28").ClearContents()
